$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "firefox" -> "chrome" (E2 keeps referencing the same shared-string slot,
# its text is simply updated in place)
$ws.Range("E2").Value = "chrome"

# New column F: header "GoogleSearch" and a new test-case value "GitHub"
$ws.Range("F1").Value = "GoogleSearch"
$ws.Range("F2").Value = "GitHub"

# Size column F to fit its new contents (adds customWidth like columns A-D)
$ws.Columns("F:F").AutoFit() | Out-Null

# Move/record the active selection on the newly added cell
$ws.Range("F2").Select() | Out-Null
